$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 8) with the "turn" command entry
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "turn"
$ws.Range("C8").Value = 2

# Update the active selection to match the new last row (C8)
$ws.Range("C8").Select()
